# Update "想去人数" (people interested) counts on the 展览 (Exhibition) sheet
# and the corresponding rows on the 全部类型 (All types) aggregate sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 500
$ws1.Range("F9").Value = 477
$ws1.Range("F13").Value = 6286
$ws1.Range("F14").Value = 211
$ws1.Range("F15").Value = 331
$ws1.Range("F16").Value = 2431
$ws1.Range("F17").Value = 129
$ws1.Range("F18").Value = 218
$ws1.Range("F20").Value = 475

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 500
$ws4.Range("F11").Value = 477
$ws4.Range("F16").Value = 6287
$ws4.Range("F18").Value = 211
$ws4.Range("F19").Value = 331
$ws4.Range("F20").Value = 2431
$ws4.Range("F21").Value = 129
$ws4.Range("F22").Value = 218
$ws4.Range("F24").Value = 475
